# The sheet "ADHD MC" had an empty, unused column AG sitting between the
# "FASz" column (AF) and the "cubost"/"TMTAz"/"TMTBz" columns (which were
# stored in AH:AJ). This edit removes that stray empty column so the three
# trailing columns shift left into AG:AI (cubost -> AG, TMTAz -> AH,
# TMTBz -> AI), and also removes a leftover scratch row (row 112) that held
# a one-off helper formula (=E107+E80) unrelated to the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty column AG - this shifts the old AH/AI/AJ ("cubost",
# "TMTAz", "TMTBz") one column to the left, becoming AG/AH/AI, and drops
# the now-empty trailing column (formerly AJ).
$ws.Columns("AG:AG").Delete()

# Remove the leftover helper row that only contained a stray formula.
$ws.Rows("112:112").Delete()

# Best-effort reproduction of the cosmetic column-width record that the
# authoring tool leaves on the very last column after a column deletion.
$ws.Columns(16384).ColumnWidth = 10.6

# Reset the view: scroll back to the top-left corner and select the first
# cell of the now-shifted "cubost" column (AG1).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AG1").Select() | Out-Null
